$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8, 20, 11, "2025-03-28 18:18:25", 15000.015, 0, "O1118"),
    @(9, 20, 11, "2025-03-28 18:41:45", 10000.01, 0, "O1141"),
    @(10, 16, 11, "2025-03-28 19:07:07", 9000.09, 0, "O1207"),
    @(11, 16, 9, "2025-03-28 19:17:33", 27000.027, 0, "O1217"),
    @(12, 16, 11, "2025-03-28 19:18:43", 10500.0105, 0, "O1218"),
    @(13, 16, 10, "2025-03-28 19:19:42", 10000.01, 0, "O1219"),
    @(14, 16, 9, "2025-03-28 19:21:47", 3500, 0, "O1221"),
    @(15, 16, 10, "2025-03-28 19:22:09", 10500.0105, 0, "O1222")
)

$startRow = 9
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
